$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()

$ws.Range("B2").Value = "Middleweight"
$ws.Range("C2").Value = "Edmen Shahbazyan"
$ws.Range("D2").Value = "Derek Brunson"

$ws.Range("B3").Value = "Women's Flyweight"
$ws.Range("C3").Value = "Jennifer Maia"
$ws.Range("D3").Value = "Joanne Calderwood"

$ws.Range("B4").Value = "Welterweight"
$ws.Range("C4").Value = "Randy Brown"
$ws.Range("D4").Value = "Vicente Luque"

$ws.Range("B5").Value = "Lightweight"
$ws.Range("C5").Value = "Justin Gaethje"
$ws.Range("D5").Value = "Khabib Nurmagomedov"

$ws.Range("E2:E5").Font.Name = "Arial"

$ws.Range("C10").Select() | Out-Null
